$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Type" column (B) and "Compass" column (H): shift data so
# Site/Lat/Long/Elevation/Slope/Aspect occupy B:G, with a new numeric id in A.

# New header row: A1 blank, B1:G1 = Site, Latitude, Longitude, Elevation, Slope, Aspect
$ws.Range("A1").Value = $null
$ws.Range("B1").Value = "Site"
$ws.Range("C1").Value = "Latitude"
$ws.Range("D1").Value = "Longitude"
$ws.Range("E1").Value = "Elevation"
$ws.Range("F1").Value = "Slope"
$ws.Range("G1").Value = "Aspect"
$ws.Range("H1").Value = $null

# Site codes moved from column A to column B; add sequential numeric id in A
$siteCodes = @("GOR", "SCT", "STS", "WON")
for ($i = 0; $i -lt 4; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $i + 1
    $ws.Cells.Item($r, 2).Value = $siteCodes[$i]
}

# Clear the old "Type" text values that used to sit in column B (now overwritten above)
# and the old "Compass" values that used to sit in column H.
$ws.Range("H2:H5").Value = $null

# Recomputed Elevation values (column E)
$ws.Range("E2").Value = 33.179824561333298
$ws.Range("E3").Value = 306.1622807
$ws.Range("E4").Value = 196.42543858666701
$ws.Range("E5").Value = 17.565789474666701

# New number format "0.0" for Elevation, Slope, Aspect columns (E:G)
$ws.Range("E2:G5").NumberFormat = "0.0"

# Match the saved selection state
[void]$ws.Range("C2:G5").Select()
